$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245, shifting existing rows 245-254 down to 246-255
$ws.Rows.Item(245).Insert()

# Populate the new row 245 with data (copy template fields from the row that is
# now at 246, then override the date / volume / price fields per the diff)
$ws.Cells.Item(245, 1).Value = 4
$ws.Cells.Item(245, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(245, 3).Value = "Los Lagos"
$ws.Cells.Item(245, 4).Value = 44747
$ws.Cells.Item(245, 5).Value = 10
$ws.Cells.Item(245, 6).Value = 100112032
$ws.Cells.Item(245, 7).Value = "Zapallo italiano"
$ws.Cells.Item(245, 8).Value = "Sin especificar"
$ws.Cells.Item(245, 9).Value = "Primera"
$ws.Cells.Item(245, 10).Value = 200
$ws.Cells.Item(245, 11).Value = 17000
$ws.Cells.Item(245, 12).Value = 18000
$ws.Cells.Item(245, 13).Value = 17500
$ws.Cells.Item(245, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(245, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(245, 16).Value = 350
$ws.Cells.Item(245, 17).Value = 50
$ws.Cells.Item(245, 18).Value = "Hortaliza"

# Match the number format / style of column D (date-like numeric format) for the new cell
$ws.Cells.Item(245, 4).NumberFormat = $ws.Cells.Item(246, 4).NumberFormat
